$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.234.63"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "1.862.29"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08181"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3036"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08173"
$ws.Range("D11").Style = "Normal"

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.173"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.801.30"
$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7075"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "29.233.23"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007887"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.786"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "2.102.39"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.395"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.955"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1449"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.953"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.424"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.95%  "

$ws.Range("E31").Value = "  -0.63%  "

$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7074"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "

$ws.Range("E38").Value = "  +0.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.142.95"
$ws.Range("E41").Value = "  +6.24%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9199"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4275"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.870"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9994"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.774"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("D49").Value = "1.997.93"
$ws.Range("E49").Value = "  -0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.215"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.956"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.38%  "

